# Regenerate merged AHB files
# - Rename the diff-comparison header columns from the old "_old"/"_new"
#   suffix convention to the versioned "_FV2404"/"_FV2410" convention.
# - Freeze the header row.
# - Turn the used range into a proper Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) ------------------------------------------
$oldHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$newHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
# Column K (11) holds "diff" and is left untouched.
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into an Excel Table --------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U79"), 0, 1)
$tbl.Name = "Table1"
